$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45970
$ws.Range("B2").Value = 72.98999999999999
$ws.Range("C2").Value = 61.78
$ws.Range("D2").Value = 55.45
$ws.Range("E2").Value = 53.36
$ws.Range("F2").Value = 53.12
$ws.Range("G2").Value = 56.27
$ws.Range("H2").Value = 67.66
$ws.Range("I2").Value = 74.31999999999999
$ws.Range("J2").Value = 67.89
$ws.Range("K2").Value = 40.43
$ws.Range("L2").Value = 2.67
$ws.Range("M2").Value = 0.65
$ws.Range("N2").Value = 0.68
$ws.Range("O2").Value = 4.31
$ws.Range("P2").Value = 4.29
$ws.Range("Q2").Value = 17.64
$ws.Range("R2").Value = 56.83
$ws.Range("S2").Value = 84.29000000000001
$ws.Range("T2").Value = 95.01000000000001
$ws.Range("U2").Value = 105.5
$ws.Range("V2").Value = 112.92
$ws.Range("W2").Value = 115.12
$ws.Range("X2").Value = 104.01
$ws.Range("Y2").Value = 89.48
$ws.Range("Z2").Value = 58.19
$ws.Range("AB2").Value = 105.38
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 114.02
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 100.26
$ws.Range("AG2").Value = "2h-16h"
